$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# New rows of "Colville" surface fuels data (radcliffe2 study), appended after row 1017
$newRows = @(
    @(6.608,"duff","Mg/ha",1,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(1.285,"fwd","Mg/ha",1,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.113,"herb","Mg/ha",1,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.693,"hundred_hour","Mg/ha",1,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(3.421,"litter","Mg/ha",1,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.381,"one_hour","Mg/ha",1,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.087,"shrub","Mg/ha",1,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.211,"ten_hour","Mg/ha",1,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(11.029,"duff","Mg/ha",5,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(1.775,"fwd","Mg/ha",5,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.54,"herb","Mg/ha",5,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.984,"hundred_hour","Mg/ha",5,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(4.267,"litter","Mg/ha",5,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.376,"one_hour","Mg/ha",5,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.202,"shrub","Mg/ha",5,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.415,"ten_hour","Mg/ha",5,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(7.673,"duff","Mg/ha",10,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(2.018,"fwd","Mg/ha",10,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.346,"herb","Mg/ha",10,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.986,"hundred_hour","Mg/ha",10,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(4.366,"litter","Mg/ha",10,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.504,"one_hour","Mg/ha",10,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.776,"shrub","Mg/ha",10,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.528,"ten_hour","Mg/ha",10,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(1.57,"fwd","Mg/ha",15,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.669,"herb","Mg/ha",15,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.8,"hundred_hour","Mg/ha",15,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.466,"one_hour","Mg/ha",15,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.55,"shrub","Mg/ha",15,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.304,"ten_hour","Mg/ha",15,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(7.533,"duff","Mg/ha",20,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(2.614,"fwd","Mg/ha",20,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.447,"herb","Mg/ha",20,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(1.431,"hundred_hour","Mg/ha",20,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(3.575,"litter","Mg/ha",20,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.558,"one_hour","Mg/ha",20,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.452,"shrub","Mg/ha",20,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA"),
    @(0.626,"ten_hour","Mg/ha",20,"thinburn","radcliffe2","both","washington","NA","commercial","NA","NA")

)

$startRow = 1018
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $row = $newRows[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

# Update the view: scroll/select to match the post-edit state
$ws.Activate() | Out-Null
$ws.Range("E1008").Select() | Out-Null

Write-Output "Added $($newRows.Length) rows to 'data' sheet."
